$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.166.75'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '2.763.45'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.56%  '
$ws.Range("E9").Value = '  -2.74%  '
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  +3.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '3.251.47'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").Value = '63.797.93'
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '2.770.26'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '361.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("E21").Value = '  -4.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.527'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.16%  '
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '0.0₃0910'
$ws.Range("E28").Value = '  -4.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.99%  '
$ws.Range("E30").Value = '  -3.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.50'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.36%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '348.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.32%  '
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '137.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("E47").Value = '  -2.85%  '
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.23%  '
